$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin prices / 1h volume deltas (GitHub Actions scrape).

$ws.Range("D2").Value = "25.887.86"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "1.638.26"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'214.60"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").Value = "'0.5021"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  -0.48%  "

$ws.Range("D8").Value = "'0.2567"
$ws.Range("E8").Value = "  -0.45%  "

$ws.Range("D9").Value = "'0.06376"
$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").Value = "'19.50"
$ws.Range("E10").Value = "  -0.12%  "

$ws.Range("D11").Value = "'0.07787"
$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("D12").Value = "1.651.44"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "'4.263"
$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("D14").Value = "1.865.62"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").Value = "'0.5407"
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").Value = "0.0₅7849"
$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("D17").Value = "'64.59"
$ws.Range("E17").Value = "  +1.59%  "

$ws.Range("D18").Value = "25.959.20"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "  -0.55%  "

$ws.Range("D20").Value = "'197.26"
$ws.Range("E20").Value = "  -3.38%  "

$ws.Range("D21").Value = "'4.370"
$ws.Range("E21").Value = "  +1.43%  "

$ws.Range("D22").Value = "'9.899"
$ws.Range("E22").Value = "  -0.95%  "

$ws.Range("D23").Value = "'5.961"
$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  -0.43%  "

$ws.Range("D25").Value = "'1.875"
$ws.Range("E25").Value = "  -3.20%  "

$ws.Range("D26").Value = "'139.49"
$ws.Range("E26").Value = "  -1.44%  "

$ws.Range("D27").Value = "'0.1140"
$ws.Range("E27").Value = "  -1.31%  "

$ws.Range("D28").Value = "'6.829"
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("D29").Value = "'15.68"
$ws.Range("E29").Value = "  -0.56%  "

$ws.Range("D30").Value = "'1.239"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").Value = "'0.04859"
$ws.Range("E31").Value = "  -4.19%  "

$ws.Range("D32").Value = "'3.256"
$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").Value = "'3.183"
$ws.Range("E33").Value = "  -0.32%  "

$ws.Range("D34").Value = "'1.526"
$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("D35").Value = "'2.362"
$ws.Range("E35").Value = "  +0.95%  "

$ws.Range("D36").Value = "'0.8869"
$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("D37").Value = "'2.605"
$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.5523"
$ws.Range("E38").Value = "  -2.14%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.129.26"
$ws.Range("E39").Value = "  -0.83%  "

$ws.Range("D40").Value = "'0.01558"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("D42").Value = "'5.678"
$ws.Range("E42").Value = "  +0.72%  "

$ws.Range("D43").Value = "'0.8142"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").Value = "'99.28"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").Value = "1.776.20"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("E46").Value = "  +6.34%  "

$ws.Range("D47").Value = "'0.4529"
$ws.Range("E47").Value = "  +0.23%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").Value = "'55.02"
$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("D50").Value = "'0.05088"
$ws.Range("E50").Value = "  +1.36%  "

$ws.Range("D51").Value = "'1.005"
$ws.Range("E51").Value = "  -0.10%  "
